# Updates cryptos list price (D) and volume(1h) (E) columns per Feb 26 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" values (column D) keyed by row number. Some rows keep their
# existing price (not present in this map) and only the volume changes.
$priceUpdates = @{
    2 = "23.488.79"
    3 = "1.633.65"
    4 = "1.002"
    5 = "308.29"
    7 = "0.3771"
    8 = "52.77"
    9 = "0.3685"
    10 = "1.272"
    11 = "0.08184"
    13 = "23.13"
    14 = "6.650"
    15 = "0.00001279"
    16 = "7.454"
    17 = "1.638.89"
    18 = "94.79"
    19 = "0.06950"
    20 = "18.38"
    21 = "6.576"
    22 = "1.004"
    23 = "23.443.51"
    24 = "12.95"
    25 = "3.113"
    27 = "21.40"
    28 = "151.36"
    29 = "5.335"
    30 = "136.40"
    31 = "2.424"
    32 = "6.780"
    33 = "1.814.49"
    34 = "0.9728"
    35 = "0.02812"
    36 = "10.46"
    37 = "0.07421"
    39 = "0.2536"
    40 = "0.08829"
    41 = "1.392"
    42 = "0.7149"
    43 = "12.59"
    44 = "16.17"
    45 = "0.6606"
    46 = "2.353"
    48 = "4.047"
    49 = "0.08038"
    50 = "131.07"
    51 = "1.215"
}

# New "Volume(1h)" values (column E) keyed by row number.
$volumeUpdates = @{
    2 = "  +1.95%  "
    3 = "  +3.10%  "
    4 = "  +0.29%  "
    5 = "  +2.64%  "
    6 = "  +0.44%  "
    7 = "  +0.34%  "
    8 = "  +3.63%  "
    9 = "  +2.70%  "
    10 = "  +3.08%  "
    11 = "  +2.27%  "
    12 = "  +0.65%  "
    13 = "  +4.32%  "
    14 = "  +2.58%  "
    15 = "  +3.63%  "
    16 = "  +1.60%  "
    17 = "  +3.40%  "
    18 = "  +2.25%  "
    19 = "  +3.08%  "
    20 = "  +3.06%  "
    21 = "  +2.97%  "
    22 = "  +0.50%  "
    23 = "  +1.80%  "
    24 = "  +1.96%  "
    25 = "  +8.97%  "
    26 = "  +1.72%  "
    27 = "  +3.09%  "
    28 = "  +1.99%  "
    29 = "  +3.29%  "
    30 = "  +3.46%  "
    31 = "  +3.13%  "
    32 = "  +2.96%  "
    33 = "  +3.29%  "
    34 = "  +3.12%  "
    35 = "  +5.51%  "
    36 = "  +4.65%  "
    37 = "  +0.14%  "
    38 = "  +2.72%  "
    39 = "  +1.89%  "
    40 = "  +0.63%  "
    41 = "  +3.74%  "
    42 = "  +2.81%  "
    43 = "  +4.24%  "
    44 = "  +9.15%  "
    45 = "  +3.11%  "
    46 = "  +3.87%  "
    47 = "  +0.47%  "
    48 = "  +1.70%  "
    49 = "  +2.24%  "
    50 = "  +0.20%  "
    51 = "  +1.41%  "
}

# Column D holds prices formatted like "23.488.79" / "1.002" which Excel would
# otherwise auto-coerce into numbers; force the whole data range to Text first
# so the values are written verbatim, then drop the explicit format again so
# the cells keep their original (default) style.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $priceUpdates[$row]
}

$priceRange.ClearFormats()

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $volumeUpdates[$row]
}
